$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B16").Value = '<h4>Das Autor:inquiz ist nun beendet.</h4> Bitte klicken Sie auf „Weiter", um den nächsten Test zu beginnen.'
